$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update student names in shared strings (referenced by B1/B2)
$ws.Range("B1").Value = "Haikal Ghiffari"
$ws.Range("B2").Value = "Sandra Agnes"

# Fix NRP / ID values
$ws.Range("C1").Value = 5111840000001
$ws.Range("C2").Value = 5111840000002

# Add missing row number for row 2
$ws.Range("A2").Value = 2

# Update active selection to F5
$ws.Range("F5").Select()
